$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fecha de legalización
$ws.Range("B6").Value = (Get-Date -Year 2023 -Month 5 -Day 1 -Hour 0 -Minute 0 -Second 0).Date

# Nombre del viajero
$ws.Range("B7").Value = "Juan Pérez"

# Centro de Costos
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "1234"

# Ciudad de destino
$ws.Range("B11").Value = "Bogotá"

# Fecha de Salida / Fecha de Regreso
$ws.Range("B12").Value = (Get-Date -Year 2023 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("F12").Value = (Get-Date -Year 2023 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0).Date

# Motivo del Viaje
$ws.Range("B13").Value = "Reunión de negocios"

# Relación de gastos - fila 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "1"
$ws.Range("C20").Value = "pepe"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "123"
$ws.Range("E20").Value = "concepto0"
$ws.Range("G20").Value = 50000

# Valor del anticipo - fila 32
$ws.Range("F32").Value = 1000000
$ws.Range("G32").Value = 0

# Firmas
$ws.Range("A40").Value = "Pedro Gómez"
$ws.Range("E40").Value = "María Rodríguez"

# Banco
$ws.Range("B45").Value = "Banco de Bogotá"
